$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inventory")

# Insert a new column before the current "Quantity" column (column D),
# shifting Quantity/Provider one column to the right.
$ws.Range("D1").EntireColumn.Insert()

# New column header and values: DiscountPrice
$ws.Range("D1").Value = "DiscountPrice"
$ws.Range("D2").Value = 1.85
$ws.Range("D3").Value = 1.5
$ws.Range("D4").Value = 999
$ws.Range("D5").Value = 5

# Update selection to match the authored workbook state
$ws.Range("E5").Select()
